$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.341.67"
$ws.Range("E2").Value = "  +2.69%  "

$ws.Range("D3").Value = "2.060.09"
$ws.Range("E3").Value = "  +5.64%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'235.85"
$ws.Range("E5").Value = "  +3.76%  "

$ws.Range("E6").Value = "  +4.81%  "

$ws.Range("D7").Value = "'57.75"
$ws.Range("E7").Value = "  +10.21%  "

$ws.Range("D9").Value = "'0.381"
$ws.Range("E9").Value = "  +5.54%  "

$ws.Range("D10").Value = "'57.58"
$ws.Range("E10").Value = "  +1.75%  "

$ws.Range("D11").Value = "'0.0758"
$ws.Range("E11").Value = "  +4.46%  "

$ws.Range("E12").Value = "  +4.91%  "

$ws.Range("D13").Value = "2.365.48"
$ws.Range("E13").Value = "  +5.75%  "

$ws.Range("D14").Value = "'14.26"
$ws.Range("E14").Value = "  +4.90%  "

$ws.Range("D15").Value = "'20.81"
$ws.Range("E15").Value = "  +7.97%  "

$ws.Range("D16").Value = "'0.773"
$ws.Range("E16").Value = "  +5.36%  "

$ws.Range("D17").Value = "'5.17"
$ws.Range("E17").Value = "  +5.23%  "

$ws.Range("D18").Value = "2.064.94"
$ws.Range("E18").Value = "  +5.78%  "

$ws.Range("D19").Value = "37.522.55"
$ws.Range("E19").Value = "  +3.40%  "

$ws.Range("D20").Value = "'6.08"
$ws.Range("E20").Value = "  +23.78%  "

$ws.Range("D21").Value = "'68.43"
$ws.Range("E21").Value = "  +3.10%  "

$ws.Range("D22").Value = "0.0₃0807"
$ws.Range("E22").Value = "  +3.26%  "

$ws.Range("D23").Value = "'224.57"
$ws.Range("E23").Value = "  +2.97%  "

$ws.Range("E24").Value = "  -0.26%  "

$ws.Range("D25").Value = "'2.44"
$ws.Range("E25").Value = "  +7.06%  "

$ws.Range("E26").Value = "  +3.42%  "

$ws.Range("D27").Value = "'162.67"
$ws.Range("E27").Value = "  +2.02%  "

$ws.Range("D28").Value = "'8.82"
$ws.Range("E28").Value = "  +5.53%  "

$ws.Range("E29").Value = "  +8.82%  "

$ws.Range("E30").Value = "  +8.98%  "

$ws.Range("D31").Value = "'19.19"
$ws.Range("E31").Value = "  +4.14%  "

$ws.Range("D32").Value = "'0.118"
$ws.Range("E32").Value = "  +3.09%  "

$ws.Range("D33").Value = "'2.64"
$ws.Range("E33").Value = "  +18.38%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0627"
$ws.Range("E34").Value = "  +5.95%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'4.45"
$ws.Range("E35").Value = "  +4.77%  "

$ws.Range("D36").Value = "'4.43"
$ws.Range("E36").Value = "  +7.44%  "

$ws.Range("E37").Value = "  +0.85%  "

$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("D39").Value = "'3.34"
$ws.Range("E39").Value = "  +6.34%  "

$ws.Range("D40").Value = "'5.82"
$ws.Range("E40").Value = "  +16.99%  "

$ws.Range("E41").Value = "  -1.54%  "

$ws.Range("D42").Value = "'4.40"
$ws.Range("E42").Value = "  +26.63%  "

$ws.Range("D43").Value = "'0.0948"
$ws.Range("E43").Value = "  +11.18%  "

$ws.Range("D44").Value = "1.467.33"
$ws.Range("E44").Value = "  +5.38%  "

$ws.Range("D45").Value = "'95.11"
$ws.Range("E45").Value = "  +11.89%  "

$ws.Range("D46").Value = "'0.0210"
$ws.Range("E46").Value = "  +6.92%  "

$ws.Range("D47").Value = "'16.09"
$ws.Range("E47").Value = "  +10.90%  "

$ws.Range("E48").Value = "  +6.23%  "

$ws.Range("E49").Value = "  +5.30%  "

$ws.Range("D50").Value = "'7.25"
$ws.Range("E50").Value = "  +9.35%  "

$ws.Range("D51").Value = "'2.92"
$ws.Range("E51").Value = "  +2.47%  "
